$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.796.15"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.931.10"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'550.11"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'131.85"
$ws.Range("E6").Value = "  +9.13%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +4.69%  "
$ws.Range("D9").Value = "2.926.81"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  +3.54%  "
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "'32.78"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "3.417.67"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  +8.27%  "
$ws.Range("D18").Value = "2.923.24"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "57.798.29"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'415.77"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'13.29"
$ws.Range("E21").Value = "  +4.26%  "
$ws.Range("E22").Value = "  +6.98%  "
$ws.Range("D23").Value = "'13.48"
$ws.Range("E23").Value = "  +8.20%  "
$ws.Range("D24").Value = "'6.97"
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("D25").Value = "'79.27"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'2.47"
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").Value = "'2.01"
$ws.Range("E29").Value = "  +6.49%  "
$ws.Range("D30").Value = "'7.37"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").Value = "'25.41"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("D33").Value = "'0.0980"
$ws.Range("E33").Value = "  +5.25%  "
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").Value = "'0.936"
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").Value = "  +13.10%  "
$ws.Range("D38").Value = "'48.28"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "'8.76"
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("E40").Value = "  +11.68%  "
$ws.Range("E41").Value = "  +3.84%  "
$ws.Range("D42").Value = "'374.72"
$ws.Range("E42").Value = "  +7.78%  "
$ws.Range("D43").Value = "'0.0344"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "2.695.69"
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'123.94"
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("D47").Value = "'0.235"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").Value = "'22.88"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("E51").Value = "  +3.45%  "
